$wb = $excel.ActiveWorkbook

# --- RAF-generation sheet: calibration value edits + selection ---
$wsGen = $wb.Worksheets.Item("RAF-generation")
$wsGen.Activate() | Out-Null
$wsGen.Range("B10").Value = 0.6
$wsGen.Range("B11").Value = 0.85
$wsGen.Range("B9").Select() | Out-Null

# --- About sheet: becomes the active/selected sheet on save ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
$wsAbout.Range("A45:A48").Select() | Out-Null
